$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header label changes ---
$ws.Range("D1").Value = "Description"
$ws.Range("F1").Value = "Quantity"
$ws.Range("G1").Value = "HS Code"
$ws.Range("H1").Value = "Country of Origin"

# --- New headers J1:N1 ---
$ws.Range("J1").Value = "Weight (ounces)"
$ws.Range("K1").Value = "Height (inches)"
$ws.Range("L1").Value = "Length (inches)"
$ws.Range("M1").Value = "Width (inches)"
$ws.Range("N1").Value = "Product Image Url"

# --- Row 2 (Champagne Bottle) ---
$ws.Range("D2").Value = "Sparkling wine (including champagne)"
$ws.Range("E2").Value = 120
$ws.Range("F2").Value = 10000
$ws.Range("G2").Value = 2204.1
$ws.Range("H2").Value = "France"
$ws.Range("I2").Value = "Warehouse-1"
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 9
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = "https://tse3.mm.bing.net/th/id/OIP.caEYlp7WHARd2U7BAZ_hXAHaHa?rs=1&pid=ImgDetMain&o=7&rm=3"

# --- Row 3 (Marine First Aid Kit) ---
$ws.Range("D3").Value = "First-aid boxes and kits"
$ws.Range("F3").Value = 20000
$ws.Range("G3").Value = 3006.5
$ws.Range("H3").Value = "Germany"
$ws.Range("I3").Value = "Warehouse-1"
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 10
$ws.Range("L3").Value = 20
$ws.Range("M3").Value = 18
$ws.Range("N3").Value = "https://tse4.mm.bing.net/th/id/OIP.oxGPo9sN1XhhTWpt6jDuvgHaE8?rs=1&pid=ImgDetMain&o=7&rm=3"
